$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.862.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.410.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.66%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.536"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.05%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0828"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.60%  "
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.785.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.426.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.769"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "40.768.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0918"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0735"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.56%  "
$ws.Range("E35").Value = "  -5.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.113"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.983.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0272"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.646.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "93.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.95%  "
